$p = $ppt.ActivePresentation
$p.ApplyTemplate("dummy.potx")
